# Insert a new weekly price record as row 19, pushing all subsequent
# rows (old 19..51) down by one (old row 51 becomes the new row 52).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(19).Insert()

$ws.Range("A19").Value = 2
$ws.Range("B19").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C19").Value = "Coquimbo"
$ws.Range("D19").Value = 44721
$ws.Range("E19").Value = 4
$ws.Range("F19").Value = 100112026
$ws.Range("G19").Value = "Haba"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 500
$ws.Range("K19").Value = 12000
$ws.Range("L19").Value = 14000
$ws.Range("M19").Value = 13000
$ws.Range("N19").Value = "$/saco 25 kilos"
$ws.Range("O19").Value = "Provincia de Limarí"
$ws.Range("P19").Value = 520
$ws.Range("Q19").Value = 25
$ws.Range("R19").Value = "Hortaliza"
